$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0.6654429999999999
$ws.Cells.Item(2,8).Value = 1.996329
$ws.Cells.Item(2,9).Value = 0.01979810471132212
$ws.Cells.Item(2,10).Value = 0.01979810471132213
$ws.Cells.Item(2,13).Value = 11.81073566666667
$ws.Cells.Item(2,14).Value = 35.432207
$ws.Cells.Item(2,15).Value = 0.3076347070004043
$ws.Cells.Item(2,16).Value = 0.3076347070004043
$ws.Cells.Item(2,17).Value = 7.859371374233665
$ws.Cells.Item(2,18).Value = 70.73434236810299
$ws.Cells.Item(2,19).Value = 0.006090584142030905
$ws.Cells.Item(2,20).Value = 0.006090584142030907

$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 0.6654429999999999
$ws.Cells.Item(3,8).Value = 1.996329
$ws.Cells.Item(3,9).Value = 0.01979810471132212
$ws.Cells.Item(3,10).Value = 0.01979810471132213
$ws.Cells.Item(3,13).Value = 12.45773566666667
$ws.Cells.Item(3,14).Value = 37.373207
$ws.Cells.Item(3,15).Value = 0.3244871420261927
$ws.Cells.Item(3,16).Value = 0.3244871420261927
$ws.Cells.Item(3,17).Value = 8.289912995233665
$ws.Cells.Item(3,18).Value = 74.60921695710299
$ws.Cells.Item(3,19).Value = 0.006424230415312216
$ws.Cells.Item(3,20).Value = 0.006424230415312217

$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 0.6654429999999999
$ws.Cells.Item(4,8).Value = 1.996329
$ws.Cells.Item(4,9).Value = 0.01979810471132212
$ws.Cells.Item(4,10).Value = 0.01979810471132213
$ws.Cells.Item(4,13).Value = 3.197710666666667
$ws.Cells.Item(4,14).Value = 9.593132
$ws.Cells.Item(4,15).Value = 0.08329089836363292
$ws.Cells.Item(4,16).Value = 0.0832908983636329
$ws.Cells.Item(4,17).Value = 2.127894179158667
$ws.Cells.Item(4,18).Value = 19.151047612428
$ws.Cells.Item(4,19).Value = 0.001649001927303293
$ws.Cells.Item(4,20).Value = 0.001649001927303293

$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.6654429999999999
$ws.Cells.Item(5,8).Value = 1.996329
$ws.Cells.Item(5,9).Value = 0.01979810471132212
$ws.Cells.Item(5,10).Value = 0.01979810471132213
$ws.Cells.Item(5,13).Value = 8.081220666666667
$ws.Cells.Item(5,14).Value = 24.243662
$ws.Cells.Item(5,15).Value = 0.2104918797744333
$ws.Cells.Item(5,16).Value = 0.2104918797744333
$ws.Cells.Item(5,17).Value = 5.377591724088666
$ws.Cells.Item(5,18).Value = 48.398325516798
$ws.Cells.Item(5,19).Value = 0.004167340276657259
$ws.Cells.Item(5,20).Value = 0.004167340276657259

$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0.6654429999999999
$ws.Cells.Item(6,8).Value = 1.996329
$ws.Cells.Item(6,9).Value = 0.01979810471132212
$ws.Cells.Item(6,10).Value = 0.01979810471132213
$ws.Cells.Item(6,13).Value = 2.844675333333333
$ws.Cells.Item(6,14).Value = 8.534026
$ws.Cells.Item(6,15).Value = 0.07409537283533685
$ws.Cells.Item(6,16).Value = 0.07409537283533686
$ws.Cells.Item(6,17).Value = 1.892969287839333
$ws.Cells.Item(6,18).Value = 17.036723590554
$ws.Cells.Item(6,19).Value = 0.001466947950018452
$ws.Cells.Item(6,20).Value = 0.001466947950018452

$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 24.43903466666667
$ws.Cells.Item(7,8).Value = 73.317104
$ws.Cells.Item(7,9).Value = 0.7271044512817749
$ws.Cells.Item(7,10).Value = 0.727104451281775
$ws.Cells.Item(7,13).Value = 11.81073566666667
$ws.Cells.Item(7,14).Value = 35.432207
$ws.Cells.Item(7,15).Value = 0.3076347070004043
$ws.Cells.Item(7,16).Value = 0.3076347070004043
$ws.Cells.Item(7,17).Value = 288.6429783965031
$ws.Cells.Item(7,18).Value = 2597.786805568528
$ws.Cells.Item(7,19).Value = 0.2236825648287586
$ws.Cells.Item(7,20).Value = 0.2236825648287586

$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 24.43903466666667
$ws.Cells.Item(8,8).Value = 73.317104
$ws.Cells.Item(8,9).Value = 0.7271044512817749
$ws.Cells.Item(8,10).Value = 0.727104451281775
$ws.Cells.Item(8,13).Value = 12.45773566666667
$ws.Cells.Item(8,14).Value = 37.373207
$ws.Cells.Item(8,15).Value = 0.3244871420261927
$ws.Cells.Item(8,16).Value = 0.3244871420261927
$ws.Cells.Item(8,17).Value = 304.4550338258364
$ws.Cells.Item(8,18).Value = 2740.095304432528
$ws.Cells.Item(8,19).Value = 0.2359360453509462
$ws.Cells.Item(8,20).Value = 0.2359360453509462

$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 24.43903466666667
$ws.Cells.Item(9,8).Value = 73.317104
$ws.Cells.Item(9,9).Value = 0.7271044512817749
$ws.Cells.Item(9,10).Value = 0.727104451281775
$ws.Cells.Item(9,13).Value = 3.197710666666667
$ws.Cells.Item(9,14).Value = 9.593132
$ws.Cells.Item(9,15).Value = 0.08329089836363292
$ws.Cells.Item(9,16).Value = 0.0832908983636329
$ws.Cells.Item(9,17).Value = 78.14896183663646
$ws.Cells.Item(9,18).Value = 703.3406565297281
$ws.Cells.Item(9,19).Value = 0.0605611829514554
$ws.Cells.Item(9,20).Value = 0.0605611829514554

$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 24.43903466666667
$ws.Cells.Item(10,8).Value = 73.317104
$ws.Cells.Item(10,9).Value = 0.7271044512817749
$ws.Cells.Item(10,10).Value = 0.727104451281775
$ws.Cells.Item(10,13).Value = 8.081220666666667
$ws.Cells.Item(10,14).Value = 24.243662
$ws.Cells.Item(10,15).Value = 0.2104918797744333
$ws.Cells.Item(10,16).Value = 0.2104918797744333
$ws.Cells.Item(10,17).Value = 197.4972320216498
$ws.Cells.Item(10,18).Value = 1777.475088194848
$ws.Cells.Item(10,19).Value = 0.1530495827426587
$ws.Cells.Item(10,20).Value = 0.1530495827426587

$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 24.43903466666667
$ws.Cells.Item(11,8).Value = 73.317104
$ws.Cells.Item(11,9).Value = 0.7271044512817749
$ws.Cells.Item(11,10).Value = 0.727104451281775
$ws.Cells.Item(11,13).Value = 2.844675333333333
$ws.Cells.Item(11,14).Value = 8.534026
$ws.Cells.Item(11,15).Value = 0.07409537283533685
$ws.Cells.Item(11,16).Value = 0.07409537283533686
$ws.Cells.Item(11,17).Value = 69.52111908674489
$ws.Cells.Item(11,18).Value = 625.6900717807041
$ws.Cells.Item(11,19).Value = 0.05387507540795613
$ws.Cells.Item(11,20).Value = 0.05387507540795615

$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 8.477506666666667
$ws.Cells.Item(12,8).Value = 25.43252
$ws.Cells.Item(12,9).Value = 0.2522207982916614
$ws.Cells.Item(12,10).Value = 0.2522207982916615
$ws.Cells.Item(12,13).Value = 11.81073566666667
$ws.Cells.Item(12,14).Value = 35.432207
$ws.Cells.Item(12,15).Value = 0.3076347070004043
$ws.Cells.Item(12,16).Value = 0.3076347070004043
$ws.Cells.Item(12,17).Value = 100.1255903524044
$ws.Cells.Item(12,18).Value = 901.13031317164
$ws.Cells.Item(12,19).Value = 0.07759187138186334
$ws.Cells.Item(12,20).Value = 0.07759187138186335

$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 8.477506666666667
$ws.Cells.Item(13,8).Value = 25.43252
$ws.Cells.Item(13,9).Value = 0.2522207982916614
$ws.Cells.Item(13,10).Value = 0.2522207982916615
$ws.Cells.Item(13,13).Value = 12.45773566666667
$ws.Cells.Item(13,14).Value = 37.373207
$ws.Cells.Item(13,15).Value = 0.3244871420261927
$ws.Cells.Item(13,16).Value = 0.3244871420261927
$ws.Cells.Item(13,17).Value = 105.6105371657378
$ws.Cells.Item(13,18).Value = 950.4948344916401
$ws.Cells.Item(13,19).Value = 0.08184240599722604
$ws.Cells.Item(13,20).Value = 0.08184240599722606

$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 8.477506666666667
$ws.Cells.Item(14,8).Value = 25.43252
$ws.Cells.Item(14,9).Value = 0.2522207982916614
$ws.Cells.Item(14,10).Value = 0.2522207982916615
$ws.Cells.Item(14,13).Value = 3.197710666666667
$ws.Cells.Item(14,14).Value = 9.593132
$ws.Cells.Item(14,15).Value = 0.08329089836363292
$ws.Cells.Item(14,16).Value = 0.0832908983636329
$ws.Cells.Item(14,17).Value = 27.10861349473778
$ws.Cells.Item(14,18).Value = 243.97752145264
$ws.Cells.Item(14,19).Value = 0.02100769687570513
$ws.Cells.Item(14,20).Value = 0.02100769687570513

$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 8.477506666666667
$ws.Cells.Item(15,8).Value = 25.43252
$ws.Cells.Item(15,9).Value = 0.2522207982916614
$ws.Cells.Item(15,10).Value = 0.2522207982916615
$ws.Cells.Item(15,13).Value = 8.081220666666667
$ws.Cells.Item(15,14).Value = 24.243662
$ws.Cells.Item(15,15).Value = 0.2104918797744333
$ws.Cells.Item(15,16).Value = 0.2104918797744333
$ws.Cells.Item(15,17).Value = 68.5086020764711
$ws.Cells.Item(15,18).Value = 616.5774186882401
$ws.Cells.Item(15,19).Value = 0.05309042995062
$ws.Cells.Item(15,20).Value = 0.05309042995062001

$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 8.477506666666667
$ws.Cells.Item(16,8).Value = 25.43252
$ws.Cells.Item(16,9).Value = 0.2522207982916614
$ws.Cells.Item(16,10).Value = 0.2522207982916615
$ws.Cells.Item(16,13).Value = 2.844675333333333
$ws.Cells.Item(16,14).Value = 8.534026
$ws.Cells.Item(16,15).Value = 0.07409537283533685
$ws.Cells.Item(16,16).Value = 0.07409537283533686
$ws.Cells.Item(16,17).Value = 24.11575410283556
$ws.Cells.Item(16,18).Value = 217.04178692552
$ws.Cells.Item(16,19).Value = 0.01868839408624694
$ws.Cells.Item(16,20).Value = 0.01868839408624695

$ws.Cells.Item(17,5).Value = 1
$ws.Cells.Item(17,6).Value = 0.3333333333333333
$ws.Cells.Item(17,7).Value = 0.02946533333333333
$ws.Cells.Item(17,8).Value = 0.088396
$ws.Cells.Item(17,9).Value = 0.0008766457152413409
$ws.Cells.Item(17,10).Value = 0.000876645715241341
$ws.Cells.Item(17,13).Value = 11.81073566666667
$ws.Cells.Item(17,14).Value = 35.432207
$ws.Cells.Item(17,15).Value = 0.3076347070004043
$ws.Cells.Item(17,16).Value = 0.3076347070004043
$ws.Cells.Item(17,17).Value = 0.3480072633302222
$ws.Cells.Item(17,18).Value = 3.132065369972
$ws.Cells.Item(17,19).Value = 0.0002696866477514297
$ws.Cells.Item(17,20).Value = 0.0002696866477514298

$ws.Cells.Item(18,5).Value = 1
$ws.Cells.Item(18,6).Value = 0.3333333333333333
$ws.Cells.Item(18,7).Value = 0.02946533333333333
$ws.Cells.Item(18,8).Value = 0.088396
$ws.Cells.Item(18,9).Value = 0.0008766457152413409
$ws.Cells.Item(18,10).Value = 0.000876645715241341
$ws.Cells.Item(18,13).Value = 12.45773566666667
$ws.Cells.Item(18,14).Value = 37.373207
$ws.Cells.Item(18,15).Value = 0.3244871420261927
$ws.Cells.Item(18,16).Value = 0.3244871420261927
$ws.Cells.Item(18,17).Value = 0.3670713339968888
$ws.Cells.Item(18,18).Value = 3.303642005972
$ws.Cells.Item(18,19).Value = 0.0002844602627081702
$ws.Cells.Item(18,20).Value = 0.0002844602627081703

$ws.Cells.Item(19,5).Value = 1
$ws.Cells.Item(19,6).Value = 0.3333333333333333
$ws.Cells.Item(19,7).Value = 0.02946533333333333
$ws.Cells.Item(19,8).Value = 0.088396
$ws.Cells.Item(19,9).Value = 0.0008766457152413409
$ws.Cells.Item(19,10).Value = 0.000876645715241341
$ws.Cells.Item(19,13).Value = 3.197710666666667
$ws.Cells.Item(19,14).Value = 9.593132
$ws.Cells.Item(19,15).Value = 0.08329089836363292
$ws.Cells.Item(19,16).Value = 0.0832908983636329
$ws.Cells.Item(19,17).Value = 0.0942216106968889
$ws.Cells.Item(19,18).Value = 0.8479944962720001
$ws.Cells.Item(19,19).Value = 0.00007301660916908081
$ws.Cells.Item(19,20).Value = 0.0000730166091690808

$ws.Cells.Item(20,5).Value = 1
$ws.Cells.Item(20,6).Value = 0.3333333333333333
$ws.Cells.Item(20,7).Value = 0.02946533333333333
$ws.Cells.Item(20,8).Value = 0.088396
$ws.Cells.Item(20,9).Value = 0.0008766457152413409
$ws.Cells.Item(20,10).Value = 0.000876645715241341
$ws.Cells.Item(20,13).Value = 8.081220666666667
$ws.Cells.Item(20,14).Value = 24.243662
$ws.Cells.Item(20,15).Value = 0.2104918797744333
$ws.Cells.Item(20,16).Value = 0.2104918797744333
$ws.Cells.Item(20,17).Value = 0.2381158606835556
$ws.Cells.Item(20,18).Value = 2.143042746152
$ws.Cells.Item(20,19).Value = 0.0001845268044973525
$ws.Cells.Item(20,20).Value = 0.0001845268044973525

$ws.Cells.Item(21,5).Value = 1
$ws.Cells.Item(21,6).Value = 0.3333333333333333
$ws.Cells.Item(21,7).Value = 0.02946533333333333
$ws.Cells.Item(21,8).Value = 0.088396
$ws.Cells.Item(21,9).Value = 0.0008766457152413409
$ws.Cells.Item(21,10).Value = 0.000876645715241341
$ws.Cells.Item(21,13).Value = 2.844675333333333
$ws.Cells.Item(21,14).Value = 8.534026
$ws.Cells.Item(21,15).Value = 0.07409537283533685
$ws.Cells.Item(21,16).Value = 0.07409537283533686
$ws.Cells.Item(21,17).Value = 0.08381930692177778
$ws.Cells.Item(21,18).Value = 0.7543737622960001
$ws.Cells.Item(21,19).Value = 0.00006495539111530769
$ws.Cells.Item(21,20).Value = 0.0000649553911153077
